# LCA_infrastructure.xlsx - "HEATING" sheet:
# Add a new reference row for "natural gas-fired boiler" (code T3 / NG),
# inserted above the existing "district heating - ..." boiler rows, and
# move the selection cursor to reflect the editor's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HEATING")

# Insert a new row at position 4; existing rows 4-6 shift down to 5-7.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the "natural gas-fired boiler" data.
$ws.Range("A4").Value = "natural gas-fired boiler"
$ws.Range("B4").Value = "T3"
$ws.Range("C4").Value = "NG"
$ws.Range("D4").Value = 0.8
$ws.Range("E4").Value = 1.403
$ws.Range("F4").Value = 0.1
$ws.Range("G4").Value = 0.22
$ws.Range("H4").Value = "KBOB 2019, costs in USD-2015"

# Update selection to match the target state.
$ws.Range("A11").Select()
